# Luukv's HMXBparameters.xlsx edit:
#  - "distance" column formula now scales parallax-derived distance by 1000
#    (1/parallax -> 1/parallax * 1000), recalculating every data row.
#  - View state: scrolled one column right (topLeftCell D1) with N6 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo  = $ws.ListObjects.Item("Tabel13")
$col = $lo.ListColumns.Item("distance")

# Update the calculated column formula cell-by-cell (keeps each row's <f> a
# plain, non-shared formula, matching how the workbook stores it).
$dataRange = $col.DataBodyRange
$rowCount = $dataRange.Rows.Count
for ($i = 1; $i -le $rowCount; $i++) {
    $cell = $dataRange.Cells.Item($i, 1)
    $cell.Formula = "=1/Tabel13[[#This Row],[parallax]] * 1000"
}

# Recalculate so cached <v> results reflect the new formula.
$excel.CalculateFull()

# Update the view: scroll right one column and select N6, like the saved
# workbook's sheetView.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N6").Select()
